# EECS 581 - Project 3 Timesheet: add a basic week of (mostly zero) hours
# for every teammate, plus Riley Meyerkorth's real entries for that week.

$wb = $excel.ActiveWorkbook

# Sheet order in the workbook (tab order 1..6):
#   1 Aiden Burke, 2 Ty Farrington, 3 Nicholas Holmes,
#   4 Riley Meyerkorth, 5 Andrew Reyes, 6 Brett Suhr

# --- Sheet 1: Aiden Burke -------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
$ws.Range("C16").Select()

# --- Sheet 2: Ty Farrington ------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Activate()
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
$ws.Range("C20").Select()

# --- Sheet 3: Nicholas Holmes ----------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Activate()
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
$ws.Range("B16:B22").Select()

# --- Sheet 4: Riley Meyerkorth ----------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Activate()
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Sprint work and sending messaegs"
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "Meeting and creating tickets"
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = "Implementing multiplayer"
$ws.Range("B22").Value = 5
$ws.Range("C22").Select()

# --- Sheet 5: Andrew Reyes ---------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Activate()
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
$ws.Range("C23").Select()

# --- Sheet 6: Brett Suhr (left active/selected, matching the saved file) ---
$ws = $wb.Worksheets.Item(6)
$ws.Activate()
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}
$ws.Range("C19").Select()
